# Update the two input cells on "Planilha1" (the fuzzy-logic inputs).
# All the other cells on the sheet are formulas driven off A3/B3 (and the
# aggregation rows 8/9/11), so changing these two values is enough for
# Excel's automatic recalculation to ripple the new results through.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A3").Value = 102
$ws.Range("B3").Value = 216

# Move the active selection from A4 to B4, matching the saved view state.
$ws.Range("B4").Select() | Out-Null
